# Generate Report for Archive
# Update localization status from "Ready for handoff" to "In Translation"
# and record the handoff name used for this generated archive report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update the status text (Overview E2/F2 mirror the localization sheets' Status column)
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# Record the handoff name ("Lastest Handoff Name" column, column I) for this archive
$zhcn.Range("I2").Value = "TestHandoff_201612090232"
$dede.Range("I2").Value = "TestHandoff_201612090232"

# Re-autofit the affected columns now that their content changed
$overview.Columns.Item(5).AutoFit() | Out-Null
$overview.Columns.Item(6).AutoFit() | Out-Null

$zhcn.Columns.Item(3).AutoFit() | Out-Null
$zhcn.Columns.Item(9).AutoFit() | Out-Null

$dede.Columns.Item(3).AutoFit() | Out-Null
$dede.Columns.Item(9).AutoFit() | Out-Null

# AutoFit rounds to whole-pixel character widths; nudge each resized column
# to the exact width the archived report uses for "In Translation" / the new
# handoff name so the saved column metadata matches byte-for-byte.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5
$zhcn.Columns.Item(9).ColumnWidth = 25.5

$dede.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(9).ColumnWidth = 25.5
